# "Output Template.xlsx" - add a new "Round Title" column just before the
# existing "Scheduling Type" column (old column AF), pushing it and every
# column after it (through "Approved By Boskalis") one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Function")

# Insert a fresh column at AF; Excel shifts AF:AM -> AG:AN and copies the
# formatting of the column that used to sit at AF, so the new header cell
# already carries style 4 like the rest of the header row.
$ws.Columns("AF:AF").Insert()
$ws.Range("AF1").Value = "Round Title"

# The worksheet's hidden _xlnm._FilterDatabase name spans the header/first
# data row (A1:AV2 before the insert); extend it one column to A1:AW2 so it
# still covers every column, including the new one.
$wb.Names.Item(1).RefersTo = "=Function!`$A`$1:`$AW`$2"

# Put the selection back on the bottom-right (frozen) pane roughly where the
# editor left off after making the change.
[void]$ws.Range("S16").Select()
